# Applies the "logic problem" data fixes to the "données11" sheet.
# Columns: A = (some ratio/metric), B = unchanged, C = (a derived total).
# Only A and C values change for the affected rows; B is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données11")

$ws.Range("A16").Value = 43.36
$ws.Range("C16").Value = 118

$ws.Range("A17").Value = 8.73
$ws.Range("C17").Value = 127

$ws.Range("A18").Value = 29.049999999999997
$ws.Range("C18").Value = 136

$ws.Range("A19").Value = 13.73
$ws.Range("C19").Value = 103

$ws.Range("A24").Value = 23.46
$ws.Range("C24").Value = 98

$ws.Range("A28").Value = 33.900000000000006
$ws.Range("C28").Value = 116

$ws.Range("A31").Value = 70.52000000000001
$ws.Range("C31").Value = 135

$ws.Range("A35").Value = 20.73
$ws.Range("C35").Value = 107

$ws.Range("A37").Value = 24.55
$ws.Range("C37").Value = 103

$ws.Range("A38").Value = 83.04
$ws.Range("C38").Value = 130

$ws.Range("A40").Value = 13.01
$ws.Range("C40").Value = 98

$ws.Range("A41").Value = 17.14
$ws.Range("C41").Value = 113

$ws.Range("A44").Value = 13.270000000000001
$ws.Range("C44").Value = 105

$ws.Range("A45").Value = 9.7199999999999989
$ws.Range("C45").Value = 113

$ws.Range("A47").Value = 5.33
$ws.Range("C47").Value = 123

$ws.Range("A52").Value = 13.059999999999999
$ws.Range("C52").Value = 86

$ws.Range("A53").Value = 19.739999999999998
$ws.Range("C53").Value = 106

$ws.Range("A55").Value = 21.490000000000002
$ws.Range("C55").Value = 130

$ws.Range("A57").Value = 58.109999999999992
$ws.Range("C57").Value = 133

$ws.Range("A63").Value = 17.990000000000002
$ws.Range("C63").Value = 90

Write-Host "Applied data corrections to rows 16-63 of données11."
